# Generate Report for handback
# Updates the "Correspond Handoff Datetime" (D2) and
# "Correspond Handback DateTime" (G2) values for the first (5d56e154...)
# entry in both the zh-cn and de-de sheets, regenerating their timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-15 15:57:57"
$wsZhCn.Range("G2").Value = "2016-01-15 15:58:45"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-15 15:58:08"
$wsDeDe.Range("G2").Value = "2016-01-15 15:59:03"
